$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.029.74"
$ws.Range("E2").Value = "  -1.63%  "

$ws.Range("D3").Value = "3.770.45"

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'621.57"
$ws.Range("E5").Value = "  +3.34%  "

$ws.Range("D6").Value = "'182.71"
$ws.Range("E6").Value = "  -0.12%  "

$ws.Range("D7").Value = "3.767.95"
$ws.Range("E7").Value = "  +2.99%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("D10").Value = "'0.169"
$ws.Range("E10").Value = "  +3.34%  "

$ws.Range("D11").Value = "'6.31"
$ws.Range("E11").Value = "  -4.92%  "

$ws.Range("D12").Value = "'0.491"
$ws.Range("E12").Value = "  -1.62%  "

$ws.Range("D13").Value = "'41.54"
$ws.Range("E13").Value = "  +1.57%  "

$ws.Range("D14").Value = "'0.0000260"
$ws.Range("E14").Value = "  +1.65%  "

$ws.Range("D15").Value = "4.384.87"
$ws.Range("E15").Value = "  +2.97%  "

$ws.Range("D16").Value = "3.751.80"
$ws.Range("E16").Value = "  +2.76%  "

$ws.Range("D17").Value = "70.066.85"
$ws.Range("E17").Value = "  -1.53%  "

$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").Value = "'7.62"
$ws.Range("E19").Value = "  +1.24%  "

$ws.Range("D20").Value = "'16.74"
$ws.Range("E20").Value = "  -1.84%  "

$ws.Range("D21").Value = "'508.68"
$ws.Range("E21").Value = "  -2.32%  "

$ws.Range("E22").Value = "  +3.97%  "

$ws.Range("D23").Value = "'0.727"
$ws.Range("E23").Value = "  -2.62%  "

$ws.Range("D24").Value = "'2.52"
$ws.Range("E24").Value = "  +1.47%  "

$ws.Range("D25").Value = "'87.21"
$ws.Range("E25").Value = "  -0.92%  "

$ws.Range("D26").Value = "'13.19"
$ws.Range("E26").Value = "  -3.05%  "

$ws.Range("D27").Value = "'11.11"
$ws.Range("E27").Value = "  +0.85%  "

$ws.Range("E28").Value = "  +20.99%  "

$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("D30").Value = "'2.53"
$ws.Range("E30").Value = "  -1.07%  "

$ws.Range("D31").Value = "'2.91"
$ws.Range("E31").Value = "  +4.24%  "

$ws.Range("D32").Value = "'7.94"
$ws.Range("E32").Value = "  -2.41%  "

$ws.Range("D33").Value = "'31.15"
$ws.Range("E33").Value = "  -2.35%  "

$ws.Range("D34").Value = "'0.115"
$ws.Range("E34").Value = "  -0.86%  "

$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("D36").Value = "'1.07"
$ws.Range("E36").Value = "  +5.36%  "

$ws.Range("D37").Value = "'6.19"
$ws.Range("E37").Value = "  +0.79%  "

$ws.Range("D38").Value = "'0.338"
$ws.Range("E38").Value = "  -2.26%  "

$ws.Range("E39").Value = "  +2.29%  "

$ws.Range("D40").Value = "'2.13"
$ws.Range("E40").Value = "  -1.93%  "

$ws.Range("D41").Value = "'50.30"
$ws.Range("E41").Value = "  -1.54%  "

$ws.Range("D42").Value = "'45.81"
$ws.Range("E42").Value = "  +1.80%  "

$ws.Range("D43").Value = "'428.12"
$ws.Range("E43").Value = "  +2.89%  "

$ws.Range("D44").Value = "'8.74"
$ws.Range("E44").Value = "  -1.22%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.85"
$ws.Range("E45").Value = "  +1.96%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "3.019.46"
$ws.Range("E46").Value = "  -3.64%  "

$ws.Range("D47").Value = "'0.0364"
$ws.Range("E47").Value = "  -1.36%  "

$ws.Range("D48").Value = "'27.63"
$ws.Range("E48").Value = "  -2.92%  "

$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("D50").Value = "'137.40"
$ws.Range("E50").Value = "  -1.59%  "

$ws.Range("D51").Value = "'2.50"
$ws.Range("E51").Value = "  +0.68%  "

